# Add 2022-Q4 data:
#  - duplicate the current "2022-Q3" sheet (so its original data is preserved
#    under the same name on a new tab)
#  - rename the original "2022-Q3" sheet to "2022-Q4" and refresh its figures
#  - insert the new 2022-Q4 summary row at the top of "总计"

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the "2022-Q3" sheet right after itself, then swap names so
#    that the duplicate keeps the old (2022-Q3) numbers and the original
#    sheet becomes the new 2022-Q4 sheet we are going to refresh.
# ------------------------------------------------------------------
$q3src = $wb.Worksheets.Item("2022-Q3")
$q3src.Copy($null, $q3src)

$dup = $wb.Worksheets.Item("2022-Q3 (2)")
$dup.Name = "2022-Q3-holder"
$q3src.Name = "2022-Q4"
$dup.Name = "2022-Q3"

# ------------------------------------------------------------------
# 2. Refresh the fund-holding figures on the (renamed) 2022-Q4 sheet.
#    Columns D-G are stored as text in this workbook, so force a text
#    number format before writing them.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2022-Q4")
$q4.Range("D2:G5").NumberFormat = "@"

$q4.Range("D2").Value = "2.55"
$q4.Range("E2").Value = "94.19"
$q4.Range("F2").Value = "7.12"
$q4.Range("G2").Value = "0.1816"
$q4.Range("H2").Value = 4

$q4.Range("D3").Value = "1.32"
$q4.Range("E3").Value = "94.19"
$q4.Range("F3").Value = "7.12"
$q4.Range("G3").Value = "0.0940"
$q4.Range("H3").Value = 4

$q4.Range("D4").Value = "1.29"
$q4.Range("E4").Value = "93.86"
$q4.Range("F4").Value = "6.08"
$q4.Range("G4").Value = "0.0784"
$q4.Range("H4").Value = 7

$q4.Range("D5").Value = "0.52"
$q4.Range("E5").Value = "93.86"
$q4.Range("F5").Value = "6.08"
$q4.Range("G5").Value = "0.0316"
$q4.Range("H5").Value = 7

# ------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift every existing row down by one
#    and write the new 2022-Q4 row at the top (row 2).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 8 is brand new - copy the index-column formatting used by the other
# rows (e.g. A7) onto A8 before filling in its value.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 4
$total.Range("D8").Value = 0.52

$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 4
$total.Range("D7").Value = 0.54

$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 5
$total.Range("D6").Value = 0.16

$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 8
$total.Range("D5").Value = 0.18

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.38

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.39
